$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.075.96'
$ws.Range("E2").Value = '  +1.35%  '

$ws.Range("D3").Value = '2.251.35'
$ws.Range("E3").Value = '  +3.38%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = '''271.24'
$ws.Range("E5").Value = '  +5.38%  '

$ws.Range("D6").Value = '''91.93'
$ws.Range("E6").Value = '  +15.02%  '

$ws.Range("D7").Value = '''0.629'
$ws.Range("E7").Value = '  +1.50%  '

$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").Value = '''0.629'
$ws.Range("E9").Value = '  +7.25%  '

$ws.Range("D10").Value = '''45.84'
$ws.Range("E10").Value = '  +7.62%  '

$ws.Range("D11").Value = '''0.0972'

$ws.Range("D12").Value = '''8.37'
$ws.Range("E12").Value = '  +21.35%  '

$ws.Range("E13").Value = '  +1.96%  '

$ws.Range("D14").Value = '2.590.51'
$ws.Range("E14").Value = '  +3.25%  '

$ws.Range("D15").Value = '''15.17'
$ws.Range("E15").Value = '  +7.15%  '

$ws.Range("D16").Value = '2.261.24'
$ws.Range("E16").Value = '  +3.63%  '

$ws.Range("D17").Value = '''0.811'
$ws.Range("E17").Value = '  +5.50%  '

$ws.Range("D18").Value = '44.043.44'
$ws.Range("E18").Value = '  +1.39%  '

$ws.Range("D19").Value = '''0.0000106'
$ws.Range("E19").Value = '  +4.09%  '

$ws.Range("E20").Value = '  +4.17%  '

$ws.Range("D21").Value = '''70.98'
$ws.Range("E21").Value = '  +1.96%  '

$ws.Range("E22").Value = '  -0.70%  '

$ws.Range("D23").Value = '''235.06'
$ws.Range("E23").Value = '  +2.61%  '

$ws.Range("E24").Value = '  +4.87%  '

$ws.Range("E25").Value = '  -0.05%  '

$ws.Range("D26").Value = '''11.50'
$ws.Range("E26").Value = '  +8.73%  '

$ws.Range("E27").Value = '  +13.95%  '

$ws.Range("E28").Value = '  +5.49%  '

$ws.Range("D29").Value = '''41.15'
$ws.Range("E29").Value = '  -2.42%  '

$ws.Range("E30").Value = '  +0.76%  '

$ws.Range("D31").Value = '''172.91'
$ws.Range("E31").Value = '  +0.02%  '

$ws.Range("E32").Value = '  +7.05%  '

$ws.Range("D33").Value = '''21.02'
$ws.Range("E33").Value = '  +3.84%  '

$ws.Range("E34").Value = '  +5.46%  '

$ws.Range("E35").Value = '  +2.17%  '

$ws.Range("E36").Value = '  +1.80%  '

$ws.Range("E37").Value = '  +1.16%  '

$ws.Range("E38").Value = '  -2.41%  '

$ws.Range("D39").Value = '''3.52'
$ws.Range("E39").Value = '  +25.79%  '

$ws.Range("D40").Value = '''13.05'
$ws.Range("E40").Value = '  +1.06%  '

$ws.Range("D41").Value = '''0.227'
$ws.Range("E41").Value = '  +15.25%  '

$ws.Range("D42").Value = '''2.20'
$ws.Range("E42").Value = '  +5.63%  '

$ws.Range("D43").Value = '''63.89'
$ws.Range("E43").Value = '  +2.77%  '

$ws.Range("E44").Value = '  +0.05%  '

$ws.Range("D45").Value = '''0.0997'
$ws.Range("E45").Value = '  +2.11%  '

$ws.Range("E46").Value = '  +2.90%  '

$ws.Range("D47").Value = '''100.44'
$ws.Range("E47").Value = '  +0.36%  '

$ws.Range("E48").Value = '  +5.46%  '

$ws.Range("E49").Value = '  +3.01%  '

$ws.Range("E50").Value = '  +1.42%  '

$ws.Range("D51").Value = '2.479.21'
$ws.Range("E51").Value = '  +3.47%  '
